$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-6 with new B/C values
$ws.Cells.Item(2, 2).Value = '001t0000009wAyPAAU'
$ws.Cells.Item(2, 3).Value = 'The Office of the Deputy Under Secretary for Health for Operations and Management (DUSHOM/10N) leads VHA operations and ensures it continues to be the benchmark for health care excellence and value through the clinical and administrative services we provide to care for Veterans and their families.The Office of the DUSHOM operates VHA health care systems, medical centers, and outpatient sites of care.  VHA’s coverage area is divided into 18 Veterans Integrated Service Networks (VISN), each a shared system of care working together to better meet local health care needs and provide Veterans greater access to care.  The DUSHOM is responsible for 24 other VHA clinical and administrative program offices that ensure that VHA program policies and regulations are executed and supported to fulfill the operating needs of VHA field operations.'
$ws.Cells.Item(3, 2).Value = '001t0000009wAyaAAE'
$ws.Cells.Item(3, 3).Value = 'Specialty Care Services (SCS/10P4D) is a large service, encompassing 16 distinct medical specialties, chaplaincy, nutrition and food services, as well as neurology centers of excellence and a SCS center for innovation.  SCS ensures the best overall clinical, preventive, spiritual, religious, and nutritional care is available to Veterans.  Both policy and program development utilize innovative approaches, technologies and interdisciplinary collaboration both within and outside of VHA promoting dignity and respect for our Veterans.SCS provides national leadership on programs and initiatives, policy matters and issues relating to care delivery, assessment of services, and outcome analysis for twenty (20) key specialty care areas.The Chief Consultant, in collaboration with Field-based National Program Directors, serves as the principal advisor to the Assistant Deputy Under Secretary for Health for Patient Care Services, Deputy Under Secretary for Health for Policy and Services, Principal Deputy Under Secretary for Health and Under Secretary for Health.'
$ws.Cells.Item(4, 2).Value = '001t0000009wAykAAE'
$ws.Cells.Item(4, 3).Value = 'Patient Care Services (PCS/10P4) provides leadership for policy and program development to enable VA to provide the best possible healthcare for our Nations’ Veterans.  PCS is dedicated to ensuring the full continuum of health care, which comprises health promotion, disease prevention, diagnostics, therapeutic and rehabilitative care, recovery and palliative care.  Utilizing innovative approaches and technologies through interdisciplinary collaboration both within and outside of VHA, PCS policy and program development supports dignity and respectful care for Veterans.  Recently, the Office of Public Health realigned with PCS expanding the mission to include post-deployment health, population health, employee health and wellness, as well as a few clinical public health programs.'
$ws.Cells.Item(5, 2).Value = 'a03t0000001i5dGAAQ'
$ws.Cells.Item(5, 3).Value = 'Initiate the deployment of data hosting, system interface development and deployment transition activities at the VISN 20 Initial Operating Capability (IOC) sites to track for Go-Live (Q2 FY2020): Mann-Grandstaff VAMC (Spokane, WA); VA Puget Sound (Seattle/Tacoma WA)'
$ws.Cells.Item(6, 2).Value = 'a7Kt0000000KyvxEAC'
$ws.Cells.Item(6, 3).Value = 'Functions of Veterans Health Administration: in general'

# Delete rows 7-10 (shifting cells up), reducing used range to A1:C6
$ws.Range("A7:C10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

